# Swap the "step 2" content between TC3 and TC4 blocks.
# TC3's step 2 lives in row 28 (columns B and D),
# TC4's step 2 lives in row 36 (columns B and D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc3_B = $ws.Range("B28").Value2
$tc3_D = $ws.Range("D28").Value2
$tc4_B = $ws.Range("B36").Value2
$tc4_D = $ws.Range("D36").Value2

$ws.Range("B28").Value2 = $tc4_B
$ws.Range("D28").Value2 = $tc4_D
$ws.Range("B36").Value2 = $tc3_B
$ws.Range("D36").Value2 = $tc3_D
